# "Generate Report for Handoff": a second file (ffff1c2dd429-320b-4a8c-9d4e-c6585923bd02.md)
# has now also gone through a handoff cycle alongside the renamed/rehandled
# e397cde0-962c-48cc-be7e-42a007e8de2d.md (was 9e604063-8613-4655-8e1d-1237b36e4ccf.md).
# This appends row 3 to all three sheets and refreshes the shared values.

$wb = $excel.ActiveWorkbook

$oldUuidHash = "9e604063-8613-4655-8e1d-1237b36e4ccf.6ec7d91b8c6dcda04a62f63769ccb01233daed97"
$newUuidHash = "e397cde0-962c-48cc-be7e-42a007e8de2d.21435b8fa64a5f80af6c70edbff9ccf9dc6a03ab"

$oldMdName = "9e604063-8613-4655-8e1d-1237b36e4ccf.md"
$newMdName = "e397cde0-962c-48cc-be7e-42a007e8de2d.md"
$newMdName2 = "ffff1c2dd429-320b-4a8c-9d4e-c6585923bd02.md"

$hyperlinkColor = 15570276  # BGR long for RGB 64 95 ED (theme "HyperLink" font color)

function Style-AsHyperlink($rng) {
    $rng.Font.Underline = 2
    $rng.Font.Color = $hyperlinkColor
}

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value2 = $newMdName
$wsOverview.Range("D2").Value2 = "2016-03-20 05:04:56"

$wsOverview.Range("A3").Value2 = $newMdName2
$wsOverview.Range("B3").Value2 = "Ready for handoff"
$wsOverview.Range("C3").Value2 = "Ready for handoff"
$wsOverview.Range("D3").Value2 = "2016-03-20 05:04:56"
$wsOverview.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/48faafa511917ab7710738ca2785669c455a5596/e2e/$newMdName", "", "", $newMdName)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/48faafa511917ab7710738ca2785669c455a5596/e2e/$newMdName2", "", "", $newMdName2)
Style-AsHyperlink($wsOverview.Range("A2"))
Style-AsHyperlink($wsOverview.Range("A3"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$zhXlf = "$newUuidHash.zh-cn.xlf"

$wsZh.Range("A2").Value2 = $newMdName
$wsZh.Range("D2").Value2 = $zhXlf
$wsZh.Range("E2").Value2 = "2016-03-20 05:04:47"

$wsZh.Range("A3").Value2 = $newMdName2
$wsZh.Range("B3").Value2 = ".md"
$wsZh.Range("C3").Value2 = "Include"
$wsZh.Range("D3").Value2 = $zhXlf
$wsZh.Range("E3").Value2 = "2016-03-20 05:04:47"
$wsZh.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("H3").Value2 = "0001-01-01 00:00:00"
$wsZh.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("J3").Value2 = "Include"

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/48faafa511917ab7710738ca2785669c455a5596/e2e/$newMdName", "", "", $newMdName)
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/da1ef09ca6937e05a007d0a42d07ee3b55f0281b/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/$zhXlf", "", "", $zhXlf)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/48faafa511917ab7710738ca2785669c455a5596/e2e/$newMdName2", "", "", $newMdName2)
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/da1ef09ca6937e05a007d0a42d07ee3b55f0281b/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/$zhXlf", "", "", $zhXlf)
Style-AsHyperlink($wsZh.Range("A2"))
Style-AsHyperlink($wsZh.Range("D2"))
Style-AsHyperlink($wsZh.Range("A3"))
Style-AsHyperlink($wsZh.Range("D3"))

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$deXlf = "$newUuidHash.de-de.xlf"

$wsDe.Range("A2").Value2 = $newMdName
$wsDe.Range("D2").Value2 = $deXlf

$wsDe.Range("A3").Value2 = $newMdName2
$wsDe.Range("B3").Value2 = ".md"
$wsDe.Range("C3").Value2 = "Include"
$wsDe.Range("D3").Value2 = $deXlf
$wsDe.Range("E3").Value2 = "0001-01-01 00:00:00"
$wsDe.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("H3").Value2 = "0001-01-01 00:00:00"
$wsDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("J3").Value2 = "Include"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/48faafa511917ab7710738ca2785669c455a5596/e2e/$newMdName", "", "", $newMdName)
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0b951771d83aa592dac287539caa2a7414ea07eb/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/$deXlf", "", "", $deXlf)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/48faafa511917ab7710738ca2785669c455a5596/e2e/$newMdName2", "", "", $newMdName2)
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0b951771d83aa592dac287539caa2a7414ea07eb/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/$deXlf", "", "", $deXlf)
Style-AsHyperlink($wsDe.Range("A2"))
Style-AsHyperlink($wsDe.Range("D2"))
Style-AsHyperlink($wsDe.Range("A3"))
Style-AsHyperlink($wsDe.Range("D3"))

Write-Host "done"
